$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C to hold the question's language.
$ws.Columns.Item(3).Insert()

# Header for the new column.
$ws.Range("C1").Value = "Language"

# Per-row language values (rows 2-4 are Bible questions, 5-7 are Commentary).
$ws.Range("C2").Value = "English"
$ws.Range("C3").Value = "French"
$ws.Range("C4").Value = "Spanish"
$ws.Range("C5").Value = "English"
$ws.Range("C6").Value = "English"
$ws.Range("C7").Value = "English"

# Match the author's final selection/cursor position.
$ws.Range("C5").Select()
